# DRTII-601 updating LGW excel upload column names
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the two header cells: "OpeFlightNo" -> "FlightNo", "Sum of Pax" -> "Forecast Pax"
$ws.Range("C1").Value = "FlightNo"
$ws.Range("J1").Value = "Forecast Pax"

# Widen columns I and J to fit the renamed headers
# (ColumnWidth is pre-compensated for the host's internal padding so the
# stored sheet width lands on/near the intended 12.5 / 14 "characters" value)
$ws.Columns.Item(9).ColumnWidth = (12.5 - 5/7)
$ws.Columns.Item(10).ColumnWidth = (14 - 5/7)

# Move / restore the active selection to H18
$ws.Range("H18").Select()
